$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains exact text formatting (avoid numeric auto-conversion)
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "28.058.93"
$ws.Range("E2").Value = "  -1.84%  "

# Row 3
$ws.Range("D3").Value = "1.898.00"
$ws.Range("E3").Value = "  -0.74%  "

# Row 4
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").Value = "312.39"
$ws.Range("E5").Value = "  -0.32%  "

# Row 6
$ws.Range("E6").Value = "  +0.16%  "

# Row 7
$ws.Range("D7").Value = "0.4993"
$ws.Range("E7").Value = "  -0.63%  "

# Row 8
$ws.Range("D8").Value = "0.3880"
$ws.Range("E8").Value = "  -1.75%  "

# Row 9
$ws.Range("D9").Value = "0.09122"
$ws.Range("E9").Value = "  -4.99%  "

# Row 10
$ws.Range("D10").Value = "1.132"
$ws.Range("E10").Value = "  -2.10%  "

# Row 11
$ws.Range("D11").Value = "41.78"
$ws.Range("E11").Value = "  +0.50%  "

# Row 12
$ws.Range("D12").Value = "6.351"
$ws.Range("E12").Value = "  -2.61%  "

# Row 13
$ws.Range("D13").Value = "20.75"
$ws.Range("E13").Value = "  -1.95%  "

# Row 14
$ws.Range("D14").Value = "1.897.32"
$ws.Range("E14").Value = "  -0.50%  "

# Row 15
$ws.Range("D15").Value = "7.271"
$ws.Range("E15").Value = "  -3.14%  "

# Row 16
$ws.Range("D16").Value = "1.001"
$ws.Range("E16").Value = "  +0.14%  "

# Row 17
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "91.93"
$ws.Range("E17").Value = "  -1.79%  "

# Row 18
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.00001104"
$ws.Range("E18").Value = "  -2.60%  "

# Row 19
$ws.Range("D19").Value = "0.06630"
$ws.Range("E19").Value = "  +0.00%  "

# Row 20
$ws.Range("D20").Value = "17.84"
$ws.Range("E20").Value = "  -0.47%  "

# Row 21
$ws.Range("E21").Value = "  +0.27%  "

# Row 22
$ws.Range("D22").Value = "6.192"
$ws.Range("E22").Value = "  -1.16%  "

# Row 23
$ws.Range("D23").Value = "28.116.62"
$ws.Range("E23").Value = "  -1.82%  "

# Row 24
$ws.Range("D24").Value = "11.44"
$ws.Range("E24").Value = "  +0.63%  "

# Row 25
$ws.Range("D25").Value = "2.311"
$ws.Range("E25").Value = "  +1.21%  "

# Row 26
$ws.Range("B26").Value = "LEO"
$ws.Range("C26").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D26").Value = "3.395"
$ws.Range("E26").Value = "  -0.12%  "

# Row 27
$ws.Range("B27").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C27").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D27").Value = "2.127.60"
$ws.Range("E27").Value = "  +0.23%  "

# Row 28
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").Value = "2.556"
$ws.Range("E28").Value = "  -7.65%  "

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "20.82"
$ws.Range("E29").Value = "  -2.88%  "

# Row 30
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "157.89"
$ws.Range("E30").Value = "  -0.99%  "

# Row 31
$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").Value = "126.06"
$ws.Range("E31").Value = "  -1.77%  "

# Row 32
$ws.Range("B32").Value = "ImmutableX"
$ws.Range("C32").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D32").Value = "1.078"
$ws.Range("E32").Value = "  -2.15%  "

# Row 33
$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D33").Value = "0.1057"
$ws.Range("E33").Value = "  -1.36%  "

# Row 34
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "5.584"
$ws.Range("E34").Value = "  -2.28%  "

# Row 35
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "3.618"
$ws.Range("E35").Value = "  -0.12%  "

# Row 36
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "9.684"
$ws.Range("E36").Value = "  -0.02%  "

# Row 37
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").Value = "0.06579"
$ws.Range("E37").Value = "  -2.96%  "

# Row 38
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").Value = "0.02403"
$ws.Range("E38").Value = "  -1.73%  "

# Row 39
$ws.Range("B39").Value = "Algorand"
$ws.Range("C39").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D39").Value = "0.2198"
$ws.Range("E39").Value = "  -0.58%  "

# Row 40
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").Value = "1.223"
$ws.Range("E40").Value = "  -2.95%  "

# Row 41
$ws.Range("D41").Value = "1.280"
$ws.Range("E41").Value = "  +7.04%  "

# Row 42
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "0.6482"
$ws.Range("E42").Value = "  +1.28%  "

# Row 43
$ws.Range("D43").Value = "4.944"
$ws.Range("E43").Value = "  -2.87%  "

# Row 44
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").Value = "11.35"
$ws.Range("E44").Value = "  -2.14%  "

# Row 45
$ws.Range("B45").Value = "Frax"
$ws.Range("C45").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D45").Value = "1.000"
$ws.Range("E45").Value = "  +0.16%  "

# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "0.6075"
$ws.Range("E46").Value = "  +0.37%  "

# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "13.24"
$ws.Range("E47").Value = "  -3.98%  "

# Row 48
$ws.Range("B48").Value = "WEMIXTOKEN"
$ws.Range("C48").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D48").Value = "1.293"
$ws.Range("E48").Value = "  +0.39%  "

# Row 49
$ws.Range("B49").Value = "PancakeSwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D49").Value = "3.676"
$ws.Range("E49").Value = "  +0.48%  "

# Row 50
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").Value = "1.996"
$ws.Range("E50").Value = "  -1.95%  "

# Row 51
$ws.Range("B51").Value = "Quant"
$ws.Range("C51").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D51").Value = "121.74"
$ws.Range("E51").Value = "  -2.14%  "
